# TORIBDA feat: added calculateReceipt method and Receipt Class
# Populate row 4 (calculateReceipt) with Do/Check/Action data, adjust row
# height and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new data for the "calculateReceipt" row (row 4).
$ws.Range("B4").Value = "10 mins"
$ws.Range("C4").Value = "9 mins 43 seconds"
$ws.Range("D4").Value = "have only created pseudocode and added method names for the sub-tasks"
$ws.Range("E4").Value = "-"

# Row 4 grows to accommodate wrapped text, matching the other data rows.
$ws.Rows.Item(4).RowHeight = 30

# Update the saved selection to match the new authoring state.
$excel.Goto($ws.Range("B5:E8"))
